$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 5: remove the "Arduino" row contents (B5, C5, D5) ---
$ws.Range("B5:D5").Clear()

# --- Row 6: STM nucleo 2x cost becomes a formula (50-25), Purchased mark unchanged ---
$ws.Cells.Item(6, 3).Formula = "=50-25"

# --- Row 10: 485 to TTL converters, cost 10 -> 55, add Purchased checkmark ---
$ws.Cells.Item(10, 3).Value = 55
$ws.Cells.Item(10, 4).Value = "✅"
$ws.Cells.Item(10, 4).HorizontalAlignment = -4108

# --- Row 11: Inclinometer, cost 40 -> 0 ---
$ws.Cells.Item(11, 3).Value = 0

# --- Row 13: Load cell, add Purchased checkmark ---
$ws.Cells.Item(13, 4).Value = "✅"
$ws.Cells.Item(13, 4).HorizontalAlignment = -4108

# --- Row 14: display, cost 30 -> 65 ---
$ws.Cells.Item(14, 3).Value = 65

# --- Row 15: new item "Signal Conditioner", cost 45, Purchased checkmark ---
$ws.Cells.Item(15, 2).Value = "Signal Conditioner"
$ws.Cells.Item(15, 3).Value = 45
$ws.Cells.Item(15, 4).Value = "✅"
$ws.Cells.Item(15, 4).HorizontalAlignment = -4108

# --- Selection / view ---
$ws.Range("B16").Select()

$win = $excel.ActiveWindow
$win.Left = 1120
$win.Top = 880
$win.Width = 34880
$win.Height = 22500
